$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Bold the "Category" heading paragraph and drop the blank
#    paragraph that used to separate it from "Category Information".
# ------------------------------------------------------------------
$categoryHeadingIndex = $null
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $trimmed = $p.Range.Text.TrimEnd("`r")
    if ($trimmed -eq "Category") {
        $categoryHeadingIndex = $idx
    }
}

$headingPara = $d.Paragraphs($categoryHeadingIndex)
$headingStart = $headingPara.Range.Start
$headingEnd = $headingPara.Range.End

$headingXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>Category</w:t></w:r></w:p>
'@

$d.Range($headingStart, $headingEnd).InsertXML($headingXml)

# Re-locate the blank paragraph right after the (now bold) heading and
# remove it, so "Category" and "Category Information" become adjacent.
$idx = 0
$blankAfterHeadingIndex = $null
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $trimmed = $p.Range.Text.TrimEnd("`r")
    if ($trimmed -eq "Category" -and $blankAfterHeadingIndex -eq $null) {
        $blankAfterHeadingIndex = $idx + 1
    }
}
$blankPara = $d.Paragraphs($blankAfterHeadingIndex)
$blankPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Append the new "Brand" section (heading, sub-heading and table)
#    right after the last table in the document, before the document's
#    trailing empty paragraph.
# ------------------------------------------------------------------
$lastTable = $d.Tables($d.Tables.Count)
$insertPos = $lastTable.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$brandXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>B</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>rand</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Brand</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> I</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>nformation</w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="a4"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="988"/>
    <w:gridCol w:w="1842"/>
    <w:gridCol w:w="5466"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="988" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>N</w:t>
        </w:r>
        <w:r>
          <w:t>o.</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1842" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>N</w:t>
        </w:r>
        <w:r>
          <w:t>ame</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5466" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>D</w:t>
        </w:r>
        <w:r>
          <w:t>escription</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="988" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>1</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1842" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>N</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>ame</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5466" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>N</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>ame</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>of</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:t>brand or company</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> – </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>must</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>be</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>unique</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="988" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1842" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Logo</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5466" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>L</w:t>
        </w:r>
        <w:r>
          <w:t>ogo of the brand or company</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="988" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>3</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1842" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>C</w:t>
        </w:r>
        <w:r>
          <w:t>ategories</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5466" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Products made by the company fall into this category</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="988" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>4</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1842" w:type="dxa"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5466" w:type="dxa"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$insertRange.InsertXML($brandXml)

Write-Output "Brand section added"
